# Apply the changes described by the commit:
#  - rename the (only) worksheet from "Sheet1" to "impact_functions"
#  - keep the workbook's 1904-date-system setting at its existing value
#    (false); this merely makes the intent explicit since the source
#    XML attribute name changed from the LibreOffice alias
#    (dateCompatibility) to the canonical OOXML name (date1904) with no
#    actual change in semantics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "impact_functions"
$wb.Date1904 = $false
